# Apply updated crypto price/volume figures (Sun Apr  2 22:39:11 UTC 2023 refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: target cell + new text value. "AsText" forces cells whose new value
# would otherwise be auto-parsed as a number (e.g. "1.005") to stay plain text,
# matching the inlineStr/shared-string cells already used throughout column D/E,
# then restores the default "Normal" style so no extra number format sticks around.
$updates = @(
    @{ Cell = "D2"; Value = '28.052.14'; AsText = $false },
    @{ Cell = "E2"; Value = '  -1.75%  '; AsText = $false },
    @{ Cell = "D3"; Value = '1.791.00'; AsText = $false },
    @{ Cell = "E3"; Value = '  -1.92%  '; AsText = $false },
    @{ Cell = "D4"; Value = '1.005'; AsText = $true },
    @{ Cell = "E4"; Value = '  +0.41%  '; AsText = $false },
    @{ Cell = "D5"; Value = '313.59'; AsText = $true },
    @{ Cell = "E5"; Value = '  -0.65%  '; AsText = $false },
    @{ Cell = "E6"; Value = '  +0.27%  '; AsText = $false },
    @{ Cell = "D7"; Value = '0.5187'; AsText = $true },
    @{ Cell = "E7"; Value = '  +1.32%  '; AsText = $false },
    @{ Cell = "D8"; Value = '0.3810'; AsText = $true },
    @{ Cell = "E8"; Value = '  -3.71%  '; AsText = $false },
    @{ Cell = "D9"; Value = '0.07808'; AsText = $true },
    @{ Cell = "E9"; Value = '  -4.78%  '; AsText = $false },
    @{ Cell = "D10"; Value = '41.31'; AsText = $true },
    @{ Cell = "E10"; Value = '  -0.99%  '; AsText = $false },
    @{ Cell = "D11"; Value = '1.093'; AsText = $true },
    @{ Cell = "E11"; Value = '  -1.85%  '; AsText = $false },
    @{ Cell = "E12"; Value = '  +0.39%  '; AsText = $false },
    @{ Cell = "D13"; Value = '6.252'; AsText = $true },
    @{ Cell = "E13"; Value = '  -1.48%  '; AsText = $false },
    @{ Cell = "D14"; Value = '20.44'; AsText = $true },
    @{ Cell = "E14"; Value = '  -3.56%  '; AsText = $false },
    @{ Cell = "D15"; Value = '1.794.14'; AsText = $false },
    @{ Cell = "E15"; Value = '  -1.90%  '; AsText = $false },
    @{ Cell = "D16"; Value = '7.261'; AsText = $true },
    @{ Cell = "E16"; Value = '  -3.78%  '; AsText = $false },
    @{ Cell = "D17"; Value = '91.98'; AsText = $true },
    @{ Cell = "E17"; Value = '  -0.99%  '; AsText = $false },
    @{ Cell = "D18"; Value = '0.00001079'; AsText = $true },
    @{ Cell = "E18"; Value = '  -4.33%  '; AsText = $false },
    @{ Cell = "D19"; Value = '0.06530'; AsText = $true },
    @{ Cell = "E19"; Value = '  -2.02%  '; AsText = $false },
    @{ Cell = "E20"; Value = '  +0.38%  '; AsText = $false },
    @{ Cell = "D21"; Value = '17.23'; AsText = $true },
    @{ Cell = "E21"; Value = '  -3.45%  '; AsText = $false },
    @{ Cell = "D22"; Value = '5.927'; AsText = $true },
    @{ Cell = "E22"; Value = '  -2.84%  '; AsText = $false },
    @{ Cell = "D23"; Value = '28.086.67'; AsText = $false },
    @{ Cell = "E23"; Value = '  -1.74%  '; AsText = $false },
    @{ Cell = "D24"; Value = '11.09'; AsText = $true },
    @{ Cell = "E24"; Value = '  -3.05%  '; AsText = $false },
    @{ Cell = "D25"; Value = '2.259'; AsText = $true },
    @{ Cell = "E25"; Value = '  -0.15%  '; AsText = $false },
    @{ Cell = "D26"; Value = '160.62'; AsText = $true },
    @{ Cell = "D27"; Value = '20.35'; AsText = $true },
    @{ Cell = "E27"; Value = '  -5.17%  '; AsText = $false },
    @{ Cell = "D28"; Value = '1.994.14'; AsText = $false },
    @{ Cell = "D29"; Value = '2.317'; AsText = $true },
    @{ Cell = "E29"; Value = '  -3.93%  '; AsText = $false },
    @{ Cell = "D30"; Value = '122.50'; AsText = $true },
    @{ Cell = "E30"; Value = '  -3.77%  '; AsText = $false },
    @{ Cell = "D31"; Value = '0.1064'; AsText = $true },
    @{ Cell = "E31"; Value = '  -2.45%  '; AsText = $false },
    @{ Cell = "D32"; Value = '1.046'; AsText = $true },
    @{ Cell = "E32"; Value = '  -5.96%  '; AsText = $false },
    @{ Cell = "D33"; Value = '3.676'; AsText = $true },
    @{ Cell = "E33"; Value = '  +0.49%  '; AsText = $false },
    @{ Cell = "D34"; Value = '5.520'; AsText = $true },
    @{ Cell = "D35"; Value = '0.07229'; AsText = $true },
    @{ Cell = "E35"; Value = '  +2.25%  '; AsText = $false },
    @{ Cell = "D36"; Value = '12.16'; AsText = $true },
    @{ Cell = "E36"; Value = '  +7.49%  '; AsText = $false },
    @{ Cell = "D37"; Value = '0.02307'; AsText = $true },
    @{ Cell = "E37"; Value = '  -2.05%  '; AsText = $false },
    @{ Cell = "D38"; Value = '8.749'; AsText = $true },
    @{ Cell = "E38"; Value = '  -0.83%  '; AsText = $false },
    @{ Cell = "D39"; Value = '0.2129'; AsText = $true },
    @{ Cell = "E39"; Value = '  -4.67%  '; AsText = $false },
    @{ Cell = "D40"; Value = '5.056'; AsText = $true },
    @{ Cell = "E40"; Value = '  -4.41%  '; AsText = $false },
    @{ Cell = "D41"; Value = '0.6118'; AsText = $true },
    @{ Cell = "E41"; Value = '  -3.38%  '; AsText = $false },
    @{ Cell = "E42"; Value = '  -1.94%  '; AsText = $false },
    @{ Cell = "D43"; Value = '1.375'; AsText = $true },
    @{ Cell = "E43"; Value = '  -1.72%  '; AsText = $false },
    @{ Cell = "D44"; Value = '13.22'; AsText = $true },
    @{ Cell = "E44"; Value = '  -2.60%  '; AsText = $false },
    @{ Cell = "D45"; Value = '3.760'; AsText = $true },
    @{ Cell = "E45"; Value = '  +0.66%  '; AsText = $false },
    @{ Cell = "D46"; Value = '0.5899'; AsText = $true },
    @{ Cell = "E46"; Value = '  -0.82%  '; AsText = $false },
    @{ Cell = "D47"; Value = '127.52'; AsText = $true },
    @{ Cell = "E47"; Value = '  +1.68%  '; AsText = $false },
    @{ Cell = "D48"; Value = '1.224'; AsText = $true },
    @{ Cell = "E48"; Value = '  +2.36%  '; AsText = $false },
    @{ Cell = "D49"; Value = '1.909'; AsText = $true },
    @{ Cell = "E49"; Value = '  -4.47%  '; AsText = $false },
    @{ Cell = "D50"; Value = '0.06733'; AsText = $true },
    @{ Cell = "E50"; Value = '  -3.01%  '; AsText = $false },
    @{ Cell = "E51"; Value = '  -2.24%  '; AsText = $false }
)

foreach ($u in $updates) {
    $c = $ws.Range($u.Cell)
    if ($u.AsText) {
        $c.NumberFormat = "@"
        $c.Value = $u.Value
        $c.Style = "Normal"
    } else {
        $c.Value = $u.Value
    }
}
